$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999464950817,
    0.99869501508327807,
    0.99525832911625989,
    0.99857402219859925,
    0.98874849071337212,
    0.96822956724667564,
    0.96679037805932366,
    0.96226064960546287,
    0.95994205543554578,
    0.95855026349696038,
    0.95695459812470585,
    0.95417555590681258,
    0.94289110098666096,
    0.93872429285802861,
    0.93613306782669214,
    0.93362675756696567,
    0.92991899644364873,
    0.92881008511345053,
    0.99764179438351241,
    0.99052477104455239,
    0.98912627826138233,
    0.9878617727425214,
    0.98822426256763929,
    0.97520413384293247,
    0.96874726030949199,
    0.96532185092248102,
    0.96151678601908941,
    0.94806607967896595,
    0.9387301505421316,
    0.93268959996801559,
    0.92503679106467018,
    0.92335749173224979,
    0.92283747941034533
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
